# Blind model names in the Evaluations sheet (column C) to prevent reviewer bias.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluations")

# Mapping of real model names -> blinded labels
$blindMap = @{
    "claude-opus-4.5" = "Model A"
    "gemini-3-pro"    = "Model B"
    "gpt-5.1"         = "Model C"
    "kimi-k2"         = "Model D"
}

# Find the last used row on the sheet (data starts on row 2; row 1 is the header)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $modelName = $cell.Value2
    if ($blindMap.ContainsKey($modelName)) {
        $cell.Value = $blindMap[$modelName]
    }
}
